# Update cell values to reflect recalculated TPM-based NATMI metrics
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 0.416828
    "H2" = 0.833656
    "I2" = 0.8138162450343883
    "J2" = 0.7445085466268955
    "M2" = 17.4294175
    "N2" = 34.858835
    "O2" = 0.4529581854295807
    "P2" = 0.3776014560521451
    "Q2" = 7.26506923769
    "R2" = 29.06027695076
    "S2" = 0.3686247296238916
    "T2" = 0.2811275112495821
    "G3" = 0.416828
    "H3" = 0.833656
    "I3" = 0.8138162450343883
    "J3" = 0.7445085466268955
    "O3" = 0.1017429801035258
    "P3" = 0.127224648983019
    "Q3" = 1.631872032956
    "R3" = 9.791232197736
    "S3" = 0.08280009002645983
    "T3" = 0.09471983850946443
    "D4" = "Inflammatory-Mac"
    "G4" = 0.416828
    "H4" = 0.833656
    "I4" = 0.8138162450343883
    "J4" = 0.7445085466268955
    "M4" = 5.397313
    "N4" = 16.191939
    "O4" = 0.1402661392829386
    "P4" = 0.1753959862028526
    "Q4" = 2.249751183164
    "R4" = 13.498507098984
    "S4" = 0.1141508627767116
    "T4" = 0.1305838107720768
    "D5" = "MuSCs"
    "G5" = 0.416828
    "H5" = 0.833656
    "I5" = 0.8138162450343883
    "J5" = 0.7445085466268955
    "M5" = 5.69137
    "N5" = 11.38274
    "O5" = 0.1479081344978025
    "P5" = 0.1233012863987851
    "Q5" = 2.37232237436
    "R5" = 9.48928949744
    "S5" = 0.1203700426270429
    "T5" = 0.09179886153398611
    "D6" = "Neutrophils"
    "G6" = 0.416828
    "H6" = 0.833656
    "I6" = 0.8138162450343883
    "J6" = 0.7445085466268955
    "M6" = 1.988496
    "N6" = 5.965488000000001
    "O6" = 0.05167731737988258
    "P6" = 0.06461997237892773
    "Q6" = 0.8288608106880001
    "R6" = 4.973164864128
    "S6" = 0.04205584038354638
    "T6" = 0.04811012171890561
    "D7" = "Resolving-Mac"
    "G7" = 0.416828
    "H7" = 0.833656
    "I7" = 0.8138162450343883
    "J7" = 0.7445085466268955
    "M7" = 4.057513666666667
    "N7" = 12.172541
    "O7" = 0.1054472433062699
    "P7" = 0.1318566499842704
    "Q7" = 1.691285306649333
    "R7" = 10.147711839896
    "S7" = 0.0858146795967361
    "T7" = 0.09816840284288045
    "E8" = 1
    "F8" = 0.3333333333333333
    "G8" = 0.032632
    "H8" = 0.097896
    "I8" = 0.0637108152714361
    "J8" = 0.08742743851251183
    "M8" = 17.4294175
    "N8" = 34.858835
    "O8" = 0.4529581854295807
    "P8" = 0.3776014560521451
    "Q8" = 0.56875675186
    "R8" = 3.41254051116
    "S8" = 0.02885833527758892
    "T8" = 0.03301272808123386
    "E9" = 1
    "F9" = 0.3333333333333333
    "G9" = 0.032632
    "H9" = 0.097896
    "I9" = 0.0637108152714361
    "J9" = 0.08742743851251183
    "O9" = 0.1017429801035258
    "P9" = 0.127224648983019
    "Q9" = 0.127753529464
    "R9" = 1.149781765176
    "S9" = 0.006482128210541129
    "T9" = 0.0111229251762388
    "D10" = "Inflammatory-Mac"
    "E10" = 1
    "F10" = 0.3333333333333333
    "G10" = 0.032632
    "H10" = 0.097896
    "I10" = 0.0637108152714361
    "J10" = 0.08742743851251183
    "M10" = 5.397313
    "N10" = 16.191939
    "O10" = 0.1402661392829386
    "P10" = 0.1753959862028526
    "Q10" = 0.176125117816
    "R10" = 1.585126060344
    "S10" = 0.008936470088692826
    "T10" = 0.01533442179909127
    "D11" = "MuSCs"
    "E11" = 1
    "F11" = 0.3333333333333333
    "G11" = 0.032632
    "H11" = 0.097896
    "I11" = 0.0637108152714361
    "J11" = 0.08742743851251183
    "M11" = 5.69137
    "N11" = 11.38274
    "O11" = 0.1479081344978025
    "P11" = 0.1233012863987851
    "Q11" = 0.18572078584
    "R11" = 1.11432471504
    "S11" = 0.009423347834132222
    "T11" = 0.0107799156351434
    "D12" = "Neutrophils"
    "E12" = 1
    "F12" = 0.3333333333333333
    "G12" = 0.032632
    "H12" = 0.097896
    "I12" = 0.0637108152714361
    "J12" = 0.08742743851251183
    "M12" = 1.988496
    "N12" = 5.965488000000001
    "O12" = 0.05167731737988258
    "P12" = 0.06461997237892773
    "Q12" = 0.06488860147200001
    "R12" = 0.583997413248
    "S12" = 0.003292404021313073
    "T12" = 0.005649558661838918
    "D13" = "Resolving-Mac"
    "E13" = 1
    "F13" = 0.3333333333333333
    "G13" = 0.032632
    "H13" = 0.097896
    "I13" = 0.0637108152714361
    "J13" = 0.08742743851251183
    "M13" = 4.057513666666667
    "N13" = 12.172541
    "O13" = 0.1054472433062699
    "P13" = 0.1318566499842704
    "Q13" = 0.1324047859706667
    "R13" = 1.191643073736
    "S13" = 0.006718129839167937
    "T13" = 0.0115278891589656
    "E14" = 1
    "F14" = 0.3333333333333333
    "G14" = 0.06272933333333333
    "H14" = 0.188188
    "I14" = 0.1224729396941756
    "J14" = 0.1680640148605926
    "M14" = 17.4294175
    "N14" = 34.858835
    "O14" = 0.4529581854295807
    "P14" = 0.3776014560521451
    "Q14" = 1.093335740163333
    "R14" = 6.56001444098
    "S14" = 0.05547512052810025
    "T14" = 0.06346121672132914
    "E15" = 1
    "F15" = 0.3333333333333333
    "G15" = 0.06272933333333333
    "H15" = 0.188188
    "I15" = 0.1224729396941756
    "J15" = 0.1680640148605926
    "O15" = 0.1017429801035258
    "P15" = 0.127224648983019
    "Q15" = 0.2455838972253334
    "R15" = 2.210255075028
    "S15" = 0.01246076186652482
    "T15" = 0.02138188529731579
    "D16" = "Inflammatory-Mac"
    "E16" = 1
    "F16" = 0.3333333333333333
    "G16" = 0.06272933333333333
    "H16" = 0.188188
    "I16" = 0.1224729396941756
    "J16" = 0.1680640148605926
    "M16" = 5.397313
    "N16" = 16.191939
    "O16" = 0.1402661392829386
    "P16" = 0.1753959862028526
    "Q16" = 0.3385698462813334
    "R16" = 3.047128616532
    "S16" = 0.01717880641753417
    "T16" = 0.02947775363168452
    "D17" = "MuSCs"
    "E17" = 1
    "F17" = 0.3333333333333333
    "G17" = 0.06272933333333333
    "H17" = 0.188188
    "I17" = 0.1224729396941756
    "J17" = 0.1680640148605926
    "M17" = 5.69137
    "N17" = 11.38274
    "O17" = 0.1479081344978025
    "P17" = 0.1233012863987851
    "Q17" = 0.3570158458533333
    "R17" = 2.14209507512
    "S17" = 0.01811474403662738
    "T17" = 0.02072250922965561
    "D18" = "Neutrophils"
    "E18" = 1
    "F18" = 0.3333333333333333
    "G18" = 0.06272933333333333
    "H18" = 0.188188
    "I18" = 0.1224729396941756
    "J18" = 0.1680640148605926
    "M18" = 1.988496
    "N18" = 5.965488000000001
    "O18" = 0.05167731737988258
    "P18" = 0.06461997237892773
    "Q18" = 0.124737028416
    "R18" = 1.122633255744
    "S18" = 0.006329072975023132
    "T18" = 0.0108602919981832
    "D19" = "Resolving-Mac"
    "E19" = 1
    "F19" = 0.3333333333333333
    "G19" = 0.06272933333333333
    "H19" = 0.188188
    "I19" = 0.1224729396941756
    "J19" = 0.1680640148605926
    "M19" = 4.057513666666667
    "N19" = 12.172541
    "O19" = 0.1054472433062699
    "P19" = 0.1318566499842704
    "Q19" = 0.2545251273008889
    "R19" = 2.290726145708
    "S19" = 0.01291443387036585
    "T19" = 0.02216035798242439
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
